$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Elaborazione" column (column C). Its content is the constant
#    "Completato controlli qualita interni e controlli spaziali" for every data
#    row, so deleting it shifts Nota (D->C) and Fonte Dati (E->D) left by one.
$ws.Columns.Item(3).Delete()

# 2. Lazio row (row 9): Ultima rilevazione becomes the literal text "17/08/2020"
#    (not a date serial), Nota gets a long explanatory note, Fonte Dati (URL)
#    is unchanged (already shifted into column D).
$ws.Range("B9").Value = '17/08/2020'
$ws.Range("C9").Value = 'Nel 2019 I dati del Lazio sono stati passati da Walter in quanto il sito dell’Idrografico Regione Lazio era fuori servizio. A fine 2019 il sito era nuovamente funzionante. Dal sito (nel 2020) sono state acquisite le serie del 2018 e sostituite a quelle passate nel 2019 da Walter. I valori forniti da Walter e quelli dell’Idrografico erano pressoche gli stessi ma non esattamente uguali. Per I valori 2019 il sito dell’Idrografico non ha ancora fornito I dati giornalieri. Tuttavia questi sono disponibili presso la sezione OpenData del Lazio che oggi (17 febbraio 2020) ha caricato anche I dati di pioggia per dicembre 2019, completando cos’ la serie del 2019. I dati OpenData 2018 sono stati confrontati con I dati 2018 acquisiti dal sito dell’Idrografico. I valori sono esattamente gli stessi.'
$ws.Range("C9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 175.35

# 3. Valle d'Aosta row (row 12): Ultima rilevazione date changes, Nota text is
#    extended with more detail.
$ws.Range("B12").Value = 43874
$ws.Range("C12").Value = 'Dati acquisiti da sito web Centro Funzionale. Nel 2019/2020 il sito del Centro Funzionale VdA ha cambiato interfaccia e formato dati. I dati sub giornalieri erano disponibili in formato semi-orario (ogni 30 minut). Da fine 2019 I dati sono disponibili solo I dati in formato orario e giornaliero. Questo non permette di ricostruire piu’ esattamente la stessa aggregazione delle serie storiche ma non esiste altra soluzione al problema.'
$ws.Range("C12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 148.5

# 4. Minor row height adjustments on rows whose wrapped text re-flowed into the
#    now-narrower layout.
$ws.Rows.Item(4).RowHeight = 35.05
$ws.Rows.Item(5).RowHeight = 79.85
$ws.Rows.Item(8).RowHeight = 141.75

# 5. Restore the selection to where the edit left off.
$ws.Range("C13").Select()
